$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing comment off of F13 before the row insert shifts things
# around, so we can re-create it at the correct post-insert location (F15).
$existingComment = $ws.Range("F13").Comment
$commentText = $existingComment.Text()
$existingComment.Delete()

# Insert two new rows above the old row 10 ("well" / plate-well section),
# pushing every row below down by two (old row 10 -> new row 12, etc.)
$ws.Rows("9:10").Insert()

# New row 9: Container / state / metadata / plate information
$ws.Range("A9").Value = "Container"
$ws.Range("B9").Value = "state"
$ws.Range("E9").Value = "metadata"
$ws.Range("F9").Value = "plate information"

# New row 10: Container / value / fileValue / source file / uploadedPlates/ABC.sdf
$ws.Range("A10").Value = "Container"
$ws.Range("C10").Value = "value"
$ws.Range("E10").Value = "fileValue"
$ws.Range("F10").Value = "source file"
$ws.Range("H10").Value = "uploadedPlates/ABC.sdf"

# Re-add the comment at its new location (old F13 -> new F15 after the
# two-row insert).
$ws.Range("F15").AddComment($commentText)

# Update the selected cell to match the edited workbook's saved selection.
$ws.Range("E11").Select()
